$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "Lastly, there will be a final c^tolumn name",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Lastly, there will be a final column name",
    2
)
